$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the existing row 6 (and below) down to row 7
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the updated weekly record
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44841
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100114007
$ws.Range("G6").Value = "Jengibre"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 18000
$ws.Range("N6").Value = "$/caja 13 kilos"
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 1385
$ws.Range("Q6").Value = 13
$ws.Range("R6").Value = "Hortaliza"

# Keep the date formatting style consistent with the other date cells in column D
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
